$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B44").Value = "801441973795329981791b8136476d51"
$ws.Range("B74").Value = "81b8198663d8342ceb3b8c0f92fab114"
$ws.Range("B89").Value = "1616b5e7f8bed5b4d7aed86321c8e87e"
$ws.Range("B99").Value = "934acdaaaa0b3be31f1a4c83585356c0"
$ws.Range("B110").Value = "391e31b1a8bd2400f63b4fbdf2ed30bd"
$ws.Range("B126").Value = "30992a194a56e3775d7bc9fa5a64bc24"
$ws.Range("B154").Value = "43d85da78f56bd9d75573aac61681246"
$ws.Range("B160").Value = "c2e8509ead66a56e0effc41c832254f6"
$ws.Range("B161").Value = "43b27c02768b9c7c3fa9e56208ca190b"
$ws.Range("B168").Value = "a1b0e2550e24d1d6623b2a13cb8c46cb"
$ws.Range("B246").Value = "a7844963b70be534ed450364d9f7d1e9"
$ws.Range("B276").Value = "a5a8399642eb3856bc0ed3d26c605c8e"
$ws.Range("B278").Value = "c471259a9ae3506bba77c0b291834b56"
$ws.Range("B293").Value = "66fae7c05456a4b684f7c16d5b50be85"
$ws.Range("B335").Value = "7243558c326a8d0e4703562fd0c150e7"
$ws.Range("B345").Value = "d1f32890b74c9e8aba42588b693f86cc"
$ws.Range("B410").Value = "c8e5ee6496752aa5375e643c806f31c9"
$ws.Range("B446").Value = "9de5a67740a3686774a6f39010a19265"
$ws.Range("B534").Value = "2ed82fcdb9b1b04ba5ba7044bfa11fc6"
$ws.Range("B553").Value = "e6e39d10005420a90d8be6f2ac9c2afb"
$ws.Range("B566").Value = "93cf8370596863b200b01bd187da9d14"
$ws.Range("B572").Value = "2829c5fc1f67e224165dc8d654e289f4"
$ws.Range("B584").Value = "4d47f7b4d0eb5996b1aebc32ac7df567"
$ws.Range("B700").Value = "cf0a52c92f73b57c3c83178f85143e6b"
$ws.Range("B756").Value = "9397a483900340432a332a438b43feee"
$ws.Range("B761").Value = "9986aac1f2a947465545084339a92eed"
$ws.Range("B768").Value = "b45c8bde2cac9396d620eb045d985164"
$ws.Range("B786").Value = "d7c32f6feaa74b68ad82f3fb3036d04e"
$ws.Range("B811").Value = "a4788cce563128e997236a350ba4b6f3"
$ws.Range("B815").Value = "7a3b54c0f3ee2ffadbb9fb3229e030d7"
$ws.Range("B816").Value = "dc3ff660a48a009b2c263afaeeb131db"
$ws.Range("B825").Value = "ee144aaf330dcd969107a5068c1f5d28"
$ws.Range("B827").Value = "b12f29376da282e56a56ae942e4a5f02"
$ws.Range("B855").Value = "ec5110340224ff40e879ea2857e85751"
$ws.Range("B862").Value = "15adcc8626573003a2667afe259f8d2e"
$ws.Range("B869").Value = "87d5f4401301379682bc0ad75b7a1ef8"
$ws.Range("B904").Value = "cd1a090fd82a983cf3eef5f74f74fdd1"
$ws.Range("B928").Value = "075dc0b3177c298bc5836ccf2890df11"
